$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "63.294.41"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "3.411.28"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue $ws.Range("D5") "568.18"
$ws.Range("E5").Value = "  +1.90%  "
Set-TextValue $ws.Range("D6") "155.54"
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.410.93"
$ws.Range("E8").Value = "  +1.66%  "
Set-TextValue $ws.Range("D9") "0.542"
$ws.Range("E9").Value = "  +2.35%  "
Set-TextValue $ws.Range("D10") "7.40"
$ws.Range("E10").Value = "  -1.33%  "
Set-TextValue $ws.Range("D11") "0.122"
$ws.Range("E11").Value = "  +3.12%  "
Set-TextValue $ws.Range("D12") "0.434"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "3.991.64"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("E14").Value = "  -3.58%  "
Set-TextValue $ws.Range("D15") "0.0000189"
$ws.Range("E15").Value = "  +4.21%  "
Set-TextValue $ws.Range("D16") "27.03"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "63.418.70"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").Value = "3.367.86"
$ws.Range("E18").Value = "  +1.28%  "
Set-TextValue $ws.Range("D19") "6.27"
$ws.Range("E19").Value = "  -3.35%  "
Set-TextValue $ws.Range("D20") "14.10"
$ws.Range("E20").Value = "  +2.22%  "
Set-TextValue $ws.Range("D21") "384.83"
$ws.Range("E21").Value = "  -1.04%  "
Set-TextValue $ws.Range("D22") "8.13"
$ws.Range("E22").Value = "  -3.80%  "
Set-TextValue $ws.Range("D23") "1.00"
$ws.Range("E23").Value = "  +0.04%  "
Set-TextValue $ws.Range("D24") "71.62"
$ws.Range("E24").Value = "  +1.54%  "
Set-TextValue $ws.Range("D25") "0.533"
$ws.Range("E25").Value = "  -1.82%  "
Set-TextValue $ws.Range("D26") "0.0000118"
$ws.Range("E26").Value = "  +20.81%  "
Set-TextValue $ws.Range("D27") "9.38"
$ws.Range("E27").Value = "  +5.84%  "
$ws.Range("E28").Value = "  -2.13%  "
$ws.Range("E29").Value = "  -0.21%  "
Set-TextValue $ws.Range("D30") "6.02"
$ws.Range("E30").Value = "  +6.61%  "
Set-TextValue $ws.Range("D31") "2.00"
$ws.Range("E31").Value = "  +0.20%  "
Set-TextValue $ws.Range("D32") "1.34"
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D33") "6.44"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D34") "23.13"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("E35").Value = "  +0.05%  "
Set-TextValue $ws.Range("D36") "6.81"
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D37") "1.46"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D38") "157.66"
$ws.Range("E38").Value = "  -1.40%  "
Set-TextValue $ws.Range("D39") "0.0765"
$ws.Range("E39").Value = "  +2.90%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D40") "1.83"
$ws.Range("E40").Value = "  -3.11%  "
$ws.Range("D41").Value = "2.884.62"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "26.83"
$ws.Range("E42").Value = "  -1.31%  "
Set-TextValue $ws.Range("D43") "0.0317"
$ws.Range("E43").Value = "  +1.90%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D44") "4.37"
$ws.Range("E44").Value = "  +0.95%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D45") "0.758"
$ws.Range("E45").Value = "  +1.31%  "
Set-TextValue $ws.Range("D46") "41.07"
$ws.Range("E46").Value = "  +0.78%  "
Set-TextValue $ws.Range("D47") "23.52"
$ws.Range("E47").Value = "  +6.29%  "
Set-TextValue $ws.Range("D48") "1.08"
$ws.Range("E48").Value = "  +2.81%  "
Set-TextValue $ws.Range("D49") "2.17"
$ws.Range("E49").Value = "  +19.36%  "
Set-TextValue $ws.Range("D50") "6.41"
$ws.Range("E50").Value = "  +1.61%  "
Set-TextValue $ws.Range("D51") "0.836"
$ws.Range("E51").Value = "  +3.30%  "
